$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 18:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1016065
$ws.Range("C4").Value = 5709
$ws.Range("D4").Value = 139691
$ws.Range("E4").Value = 819251
$ws.Range("G4").Value = 326
$ws.Range("H4").Value = 57123

# Row 14
$ws.Range("B14").Value = 68188
$ws.Range("C14").Value = 1687
$ws.Range("E14").Value = 32372
$ws.Range("G14").Value = 131
$ws.Range("H14").Value = 4674

# Row 19
$ws.Range("E19").Value = 5365
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 1699

# Row 55
$ws.Range("B55").Value = 4252
$ws.Range("C55").Value = 132
$ws.Range("D55").Value = 778
$ws.Range("E55").Value = 3309
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 165

# Row 57
$ws.Range("B57").Value = 3741
$ws.Range("C57").Value = 12
$ws.Range("E57").Value = 529
$ws.Range("F57").Value = 19
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 89

# Rows 69/70 - Irak moves above Armenia in the ranking (sorted by total cases).
# Row 69 keeps the 2nd-rank position but is now "Irak" with updated figures;
# Row 70 becomes "Armenia" keeping its previous (unchanged) figures.
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 1928
$ws.Range("C69").Value = 81
$ws.Range("D69").Value = 1319
$ws.Range("E69").Value = 519
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 90

$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 1867
$ws.Range("C70").Value = 59
$ws.Range("D70").Value = 866
$ws.Range("E70").Value = 971
$ws.Range("F70").Value = 10
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 30
